$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: set G4 (E4 stays a formula)
$ws.Range("G4").Value = 40

# Row 6: E6 becomes a literal value, G6 filled in
$ws.Range("E6").Value = 11
$ws.Range("G6").Value = 50

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 39

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 42

# Row 14
$ws.Range("E14").Value = 10
$ws.Range("G14").Value = 45

# Row 15: D, E, F all become/stay literal values, G filled in
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 15
$ws.Range("G15").Value = 40

# Row 20
$ws.Range("E20").Value = 10
$ws.Range("G20").Value = 50

# Row 21
$ws.Range("E21").Value = 10
$ws.Range("G21").Value = 50

# Row 25
$ws.Range("E25").Value = 10
$ws.Range("G25").Value = 50

# Row 26
$ws.Range("E26").Value = 10
$ws.Range("G26").Value = 50

# Row 29
$ws.Range("E29").Value = 1
$ws.Range("G29").Value = 40

# Row 32
$ws.Range("E32").Value = 5
$ws.Range("G32").Value = 45

# Row 33
$ws.Range("E33").Value = 11
$ws.Range("G33").Value = 50

# Row 34: D34 becomes a literal value, E34 stays a formula, G34 filled in
$ws.Range("D34").Value = 10
$ws.Range("G34").Value = 50

# Row 35
$ws.Range("E35").Value = 20
$ws.Range("G35").Value = 50

# Update view: scroll/freeze back to top and select E34 as the active cell
$ws.Range("E34").Select()
